# Added option to run omp or pthread using nonblocked and blocked.
# On the "omp" sheet: rename the header label in A2 from "betalgeuse" to
# "ale", and populate the previously-empty numeric columns (B:E) for the
# class1..class7 rows (rows 3-10) with the new benchmark results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("omp")

# Row 2 header label changes to "ale" (values in B2:E2 stay as-is).
$ws.Range("A2").Value = "ale"

# New data for rows 3 through 10 (columns B, C, D, E).
$data = @(
    @(3, 40095, 57939, 138641, 11762233),
    @(4, 31662, 93823, 116276, 4029743),
    @(5, 46922, 101927, 128558, 15613955),
    @(6, 38217, 71346, 90176, 4111614),
    @(7, 375485, 495794, 434452, 4004177),
    @(8, 4878322, 5484645, 3629706, 7445210),
    @(9, 47260477, 47405558, 39485704, 35810069),
    @(10, 358520148, 467047849, 427695019, 381688955)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Update the active selection to C11, matching the saved view state.
$ws.Range("C11").Select() | Out-Null
